$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert two new rows at position 5, pushing the existing rows 5-10 down to 7-12.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# New row 5: HUESCA GARCIA, ALDAIR OMAR, group 3ARHM
$ws.Cells.Item(5,1).Value = 20330051920168
$ws.Cells.Item(5,2).Value = "HUESCA"
$ws.Cells.Item(5,3).Value = "GARCIA"
$ws.Cells.Item(5,4).Value = "ALDAIR OMAR"
$ws.Cells.Item(5,5).Value = "INGLÉS III"
$ws.Cells.Item(5,6).Value = "3ARHM"
$ws.Cells.Item(5,7).Value = 6

# New row 6: VAZQUEZ VICTORIANO, MARIAN, group 3ARHM
$ws.Cells.Item(6,1).Value = 20330051920184
$ws.Cells.Item(6,2).Value = "VAZQUEZ"
$ws.Cells.Item(6,3).Value = "VICTORIANO"
$ws.Cells.Item(6,4).Value = "MARIAN"
$ws.Cells.Item(6,5).Value = "INGLÉS III"
$ws.Cells.Item(6,6).Value = "3ARHM"
$ws.Cells.Item(6,7).Value = 6
